# Groupings for the OBS Questions workbook — remove placeholder "cruft"
# left over from the 2015->2017 question-number mapping so the file can be
# consumed by ETL scripts.

$wb = $excel.ActiveWorkbook

$ws2017 = $wb.Worksheets.Item("QuestionsGroups 2017")
$ws2015 = $wb.Worksheets.Item("QuestionsGroups 2015")

# The "Audit Report" row used placeholder tokens (t3pbs, t3ebp, ...) referring
# to columns on the 2015 sheet while the 2017 question numbers were still
# unknown. Those numbers are now known, so replace the placeholder text with
# the real question-number range.
$ws2017.Range("C4").Value = "1-102, 143-149"

# Remove the now-obsolete scratch note + 2015-to-2017 token lookup table that
# explained/resolved those placeholders (kept the cell formatting, just drop
# the leftover content).
[void]$ws2017.Range("A28:B36").ClearContents()

# Leave the selection/active sheet the way the author left it when saving.
[void]$ws2017.Activate()
[void]$ws2017.Range("C4").Select()
[void]$ws2015.Activate()
[void]$ws2015.Range("C4").Select()
